{"js": "// removed 5V from barel to pi and pico\n// 1) \"JST connectoren en \" -> \"JST connectoren en/of \" in the T3.1 connector-type description.\n// 2) Update MoSCoW priority letters for the power-related rows (T4.1, T4.2, T4.4 -> W;\n//    T5.2, T5.3 -> C).\n\nconst body = context.document.body;\n\n// --- 1) JST connectoren en -> JST connectoren en/of -------------------------------\nconst jstResults = body.search(\"JST connectoren en \", { matchCase: true });\njstResults.load(\"items\");\nawait context.sync();\n\nif (jstResults.items.length > 0) {\n  jstResults.items[0].insertText(\"JST connectoren en/of \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2) MoSCoW column updates ------------------------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table1 = tables.items[0]; // rows T1..T4.1\nconst table2 = tables.items[1]; // rows T4.2..T6.1\n\n// Table 1, row 16 (0-indexed 15), column 2 (0-indexed 1) -> T4.1 (S -> W).\nconst t41 = table1.getCell(15, 1).body.paragraphs.getFirst().getRange();\nt41.insertText(\"W\", Word.InsertLocation.replace);\n\n// Table 2, row 1 (0-indexed 0) -> T4.2 (C -> W).\nconst t42 = table2.getCell(0, 1).body.paragraphs.getFirst().getRange();\nt42.insertText(\"W\", Word.InsertLocation.replace);\n\n// Table 2, row 3 (0-indexed 2) -> T4.4 (S -> W).\nconst t44 = table2.getCell(2, 1).body.paragraphs.getFirst().getRange();\nt44.insertText(\"W\", Word.InsertLocation.replace);\n\n// Table 2, row 6 (0-indexed 5) -> T5.2 (S -> C).\nconst t52 = table2.getCell(5, 1).body.paragraphs.getFirst().getRange();\nt52.insertText(\"C\", Word.InsertLocation.replace);\n\n// Table 2, row 7 (0-indexed 6) -> T5.3 (S -> C).\nconst t53 = table2.getCell(6, 1).body.paragraphs.getFirst().getRange();\nt53.insertText(\"C\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# removed 5V from barel to pi and pico\n# 1) \"JST connectoren en \" -> \"JST connectoren/of \" in the T3.1 connector-type description.\n# 2) Update MoSCoW priority letters for the power-related rows (T4.1, T4.2, T4.4 -> W;\n#    T5.2, T5.3 -> C).\n\n$d = $word.ActiveDocument\n\n# --- 1) JST connectoren en -> JST connectoren en/of ----------------------------------\n$range = $d.Content\n$found = $range.Find.Execute(\"JST connectoren en \")\nif ($found) {\n    $range.Text = \"JST connectoren en/of \"\n}\n\n# --- 2) MoSCoW column updates ---------------------------------------------------------\n# Table 1 (spans T1..T4.1): row 16, column 2 is T4.1's priority cell (S -> W).\n$table1 = $d.Tables.Item(1)\n$table1.Cell(16, 2).Range.Text = \"W\"\n\n# Table 2 (spans T4.2..T6.1): row 1 = T4.2 (C -> W), row 3 = T4.4 (S -> W),\n# row 6 = T5.2 (S -> C), row 7 = T5.3 (S -> C).\n$table2 = $d.Tables.Item(2)\n$table2.Cell(1, 2).Range.Text = \"W\"\n$table2.Cell(3, 2).Range.Text = \"W\"\n$table2.Cell(6, 2).Range.Text = \"C\"\n$table2.Cell(7, 2).Range.Text = \"C\"\n"}
